# Switch license from BY-NC to BY-SA
# (units/4/lessons/7/resources/petascale-lesson-4.7-slides.pptx)

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 2 ("Except where otherwise noted, this work by ... ") holds
# the license statement that needs to move from CC BY-NC 4.0 to
# CC BY-SA 4.0, both in the human readable text and in the URL.
# -----------------------------------------------------------------
$licenseSlide = $p.Slides.Item(2)
$licenseShape = $licenseSlide.Shapes.Item(1)
$tr = $licenseShape.TextFrame.TextRange

$fullText = $tr.Text

# --- 1. "CC BY-NC 4.0. ..." -> "CC BY-SA 4.0. ..." -----------------
$ccIdx = $fullText.IndexOf("CC BY-NC 4.0. To view a copy of this license, visit ")
if ($ccIdx -ge 0) {
    # "BY-NC " sits right after "CC " inside that run - replacing just
    # that word (keeping the trailing space) reproduces the same run
    # split PowerPoint itself performs when a user retypes a word.
    $byNcStart = $ccIdx + ("CC ").Length + 1   # 1-based TextRange index
    $byNcLen = ("BY-NC ").Length
    $seg = $tr.Characters($byNcStart, $byNcLen)
    if ($seg.Text -eq "BY-NC ") {
        $seg.Text = "BY-SA "
    }
}

# --- 2. hyperlink text "https://creativecommons.org/licenses/by-nc/4.0"
#        -> ".../by-sa/4.0" (hyperlink target itself stays as-is) -----
$fullText = $tr.Text
$urlIdx = $fullText.IndexOf("https://creativecommons.org/licenses/by-nc/4.0")
if ($urlIdx -ge 0) {
    $tailStart = $urlIdx + ("https://").Length + 1   # 1-based TextRange index
    $tailLen = ("creativecommons.org/licenses/by-nc/4.0").Length
    $seg2 = $tr.Characters($tailStart, $tailLen)
    if ($seg2.Text -eq "creativecommons.org/licenses/by-nc/4.0") {
        $seg2.Text = "creativecommons.org/licenses/by-sa/4.0"
    }
}

# -----------------------------------------------------------------
# Tiny re-flow nudge of the textbox that holds the license text
# (position shifts by 2 EMU as part of the same edit).
# -----------------------------------------------------------------
$licenseShape.Left = 566057 / 12700.0
